$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $cellRef, $text)
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.ClearFormats()
}

Set-TextValue $ws "D2" "248.74"

Set-TextValue $ws "D3" "22.57"

Set-TextValue $ws "D4" "5.277"

Set-TextValue $ws "D5" "0.05688"

Set-TextValue $ws "D6" "3.407"

Set-TextValue $ws "D7" "6.342"

Set-TextValue $ws "D8" "0.8055"

Set-TextValue $ws "D9" "0.8963"

Set-TextValue $ws "D10" "0.1401"

Set-TextValue $ws "D11" "0.07439"

Set-TextValue $ws "D12" "0.03101"

Set-TextValue $ws "D14" "0.09376"

Set-TextValue $ws "D15" "3.875"

Set-TextValue $ws "D16" "0.001585"

Set-TextValue $ws "D17" "0.04767"

Set-TextValue $ws "B18" "UpBots"
Set-TextValue $ws "C18" "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-TextValue $ws "D18" "0.01829"
Set-TextValue $ws "E18" "17UpBotsUBXTBestin24h"

Set-TextValue $ws "B19" "One"
Set-TextValue $ws "C19" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws "D19" "0.0005811"
Set-TextValue $ws "E19" "18OneONEWorstin24h"

Set-TextValue $ws "B20" "TigerCash"
Set-TextValue $ws "C20" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws "D20" "0.006432"
Set-TextValue $ws "E20" "19TigerCashTCH"

Set-TextValue $ws "D21" "0.004988"

Set-TextValue $ws "D22" "0.0010000"

Set-TextValue $ws "D23" "0.0001501"

Set-TextValue $ws "D24" "3.696"

Set-TextValue $ws "D25" "2.201"

Set-TextValue $ws "D26" "0.3260"

Set-TextValue $ws "D27" "0.1306"

Set-TextValue $ws "D40" "0.03972"

Set-TextValue $ws "B41" "KickToken"
Set-TextValue $ws "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws "D41" "0.006736"
Set-TextValue $ws "E41" "40KickTokenKICK"

Set-TextValue $ws "B42" "BKEXToken"
Set-TextValue $ws "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D42" "0.1068"
Set-TextValue $ws "E42" "41BKEXTokenBKK"

Set-TextValue $ws "B43" "CEJI"
Set-TextValue $ws "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D43" "0.002769"
Set-TextValue $ws "E43" "42CEJICEJI"

Set-TextValue $ws "D44" "0.007721"

Set-TextValue $ws "D45" "0.00005593"

Set-TextValue $ws "D47" "0.4991"

Set-TextValue $ws "D48" "0.2057"
